$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in previously-empty rows 43-44 (new tag/question pairs under "สาม")
$ws.Range("A43").Value = "สาม"
$ws.Range("B43").Value = "จดทะเบียน"

$ws.Range("A44").Value = "สาม"
$ws.Range("B44").Value = "วิธีจดทะเบียน"

# Row 75: change the existing question text (tag "สี่" unchanged)
$ws.Range("B75").Value = "ขั้นตอนก่อนการเริ่มต้นขายของ"

# Fill in previously-empty row 83
$ws.Range("A83").Value = "สี่"
$ws.Range("B83").Value = "ขั้นตอนก่อนการเริ่มขาย"

# Fill in previously-empty row 176 - B176 must be text "10", not a number
$ws.Range("A176").Value = "สิบ"
$ws.Range("B176").NumberFormat = "@"
$ws.Range("B176").Value = "10"

# Remove the three "easter" rows (200-202) entirely
$ws.Range("A200:B202").ClearContents()

# Fill in previously-empty row 212
$ws.Range("A212").Value = "ชมบอท"
$ws.Range("B212").Value = "เก่ง"

# Fill in previously-empty row 216
$ws.Range("A216").Value = "ทำอะไร"
$ws.Range("B216").Value = "ทำ"

# Fill in previously-empty rows 220-221
$ws.Range("A220").Value = "หัวข้อ"
$ws.Range("B220").Value = "หัวข้ออะไรบ้าง"

$ws.Range("A221").Value = "หัวข้อ"
$ws.Range("B221").Value = "หัวข้อ"
